# Update the "cryptos" price/volume table with refreshed values from the
# latest GitHub Actions data pull. A few coins also changed rank position
# (their whole row - Coin/Link/Price/Volume - moved), so those rows get
# all four columns rewritten instead of just Price/Volume.
#
# Numeric-looking Price values are written with a leading apostrophe so
# Excel stores them as text (matching the sheet's existing convention of
# keeping "65.487.15"-style thousand-grouped prices, and ambiguous ones
# like "0.999", as text rather than auto-converting them to numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.487.15"
$ws.Range("E2").Value = "  -1.19%  "

$ws.Range("D3").Value = "3.437.51"
$ws.Range("E3").Value = "  -4.11%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'594.72"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("D6").Value = "'135.29"
$ws.Range("E6").Value = "  -8.40%  "

$ws.Range("D7").Value = "3.437.50"
$ws.Range("E7").Value = "  -4.08%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("D10").Value = "'7.51"
$ws.Range("E10").Value = "  -3.93%  "

$ws.Range("E11").Value = "  -9.82%  "

$ws.Range("D12").Value = "'0.378"
$ws.Range("E12").Value = "  -8.47%  "

$ws.Range("D13").Value = "4.015.34"
$ws.Range("E13").Value = "  -4.24%  "

$ws.Range("E14").Value = "  -12.15%  "

$ws.Range("D15").Value = "'26.48"
$ws.Range("E15").Value = "  -10.28%  "

$ws.Range("D16").Value = "65.378.02"
$ws.Range("E16").Value = "  -1.47%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.114"
$ws.Range("E17").Value = "  -2.25%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.398.03"
$ws.Range("E18").Value = "  -5.26%  "

$ws.Range("D19").Value = "'9.97"
$ws.Range("E19").Value = "  -9.80%  "

$ws.Range("D20").Value = "'5.75"
$ws.Range("E20").Value = "  -9.04%  "

$ws.Range("D21").Value = "'13.73"
$ws.Range("E21").Value = "  -7.41%  "

$ws.Range("D22").Value = "'392.11"
$ws.Range("E22").Value = "  -7.34%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.545"
$ws.Range("E23").Value = "  -10.59%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'73.13"
$ws.Range("E24").Value = "  -6.73%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").Value = "3.577.45"
$ws.Range("E26").Value = "  -4.10%  "

$ws.Range("D27").Value = "'0.0000106"
$ws.Range("E27").Value = "  -11.96%  "

$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "'7.29"
$ws.Range("E29").Value = "  -11.07%  "

$ws.Range("D30").Value = "'2.26"
$ws.Range("E30").Value = "  -9.20%  "

$ws.Range("D31").Value = "'8.17"
$ws.Range("E31").Value = "  -12.55%  "

$ws.Range("D32").Value = "3.442.85"
$ws.Range("E32").Value = "  -3.90%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").Value = "'0.145"
$ws.Range("E34").Value = "  -7.43%  "

$ws.Range("D35").Value = "'22.68"
$ws.Range("E35").Value = "  -9.38%  "

$ws.Range("D36").Value = "'172.00"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("D37").Value = "'1.23"
$ws.Range("E37").Value = "  -13.86%  "

$ws.Range("D38").Value = "'6.85"
$ws.Range("E38").Value = "  -11.49%  "

$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "  -8.65%  "

$ws.Range("D40").Value = "'4.82"
$ws.Range("E40").Value = "  -13.40%  "

$ws.Range("D41").Value = "'0.0774"
$ws.Range("E41").Value = "  -9.21%  "

$ws.Range("D42").Value = "'0.812"
$ws.Range("E42").Value = "  -7.73%  "

$ws.Range("D43").Value = "'43.48"
$ws.Range("E43").Value = "  -5.20%  "

$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "'4.40"
$ws.Range("E45").Value = "  -14.95%  "

$ws.Range("D46").Value = "'1.62"
$ws.Range("E46").Value = "  -12.35%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'22.78"
$ws.Range("E47").Value = "  -3.11%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D49").Value = "'6.52"
$ws.Range("E49").Value = "  -8.59%  "

$ws.Range("E50").Value = "  -15.98%  "

$ws.Range("D51").Value = "2.188.13"
$ws.Range("E51").Value = "  -8.30%  "
